$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-14 Tuesday", "2024-05-15 Wednesday"),
    @("465×4=", "526×5="),
    @("405×5=", "676×6="),
    @("151×7=", "880×6="),
    @("541×3=", "311×8="),
    @("446×5=", "930×5="),
    @("102×4=", "975×5="),
    @("519×3=", "295×9="),
    @("415×6=", "193×2="),
    @("304×7=", "555×7="),
    @("403×5=", "593×5="),
    @("679×9=", "524×9="),
    @("612×6=", "334×6="),
    @("106×6=", "263×2="),
    @("565×5=", "356×4="),
    @("120×2=", "208×7="),
    @("359×6=", "170×3="),
    @("588×4=", "400×2="),
    @("909×9=", "959×7="),
    @("794×6=", "991×6="),
    @("531×7=", "998×7="),
    @("830×6=", "114×8="),
    @("150×8=", "760×5="),
    @("838×8=", "952×5="),
    @("866×5=", "801×7="),
    @("803×2=", "404×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
